# Update classification-report metrics with the new train_df results
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.8679245283018868
$ws.Range("C2").Value = 0.9387755102040817
$ws.Range("D2").Value = 0.9019607843137256
$ws.Range("E2").Value = 49

$ws.Range("B3").Value = 0.925
$ws.Range("C3").Value = 0.8409090909090909
$ws.Range("D3").Value = 0.8809523809523809
$ws.Range("E3").Value = 44

$ws.Range("B4").Value = 0.8924731182795699
$ws.Range("C4").Value = 0.8924731182795699
$ws.Range("D4").Value = 0.8924731182795699
$ws.Range("E4").Value = 0.8924731182795699

$ws.Range("B5").Value = 0.8964622641509434
$ws.Range("C5").Value = 0.8898423005565863
$ws.Range("D5").Value = 0.8914565826330533
$ws.Range("E5").Value = 93

$ws.Range("B6").Value = 0.8949279772773382
$ws.Range("C6").Value = 0.8924731182795699
$ws.Range("D6").Value = 0.8920213246588958
$ws.Range("E6").Value = 93

$ws.Range("B7").Value = 0.7777777777777778
$ws.Range("C7").Value = 0.7142857142857143
$ws.Range("D7").Value = 0.7446808510638298
$ws.Range("E7").Value = 49

$ws.Range("B8").Value = 0.7083333333333334
$ws.Range("C8").Value = 0.7727272727272727
$ws.Range("D8").Value = 0.7391304347826088
$ws.Range("E8").Value = 44

$ws.Range("B9").Value = 0.7419354838709677
$ws.Range("C9").Value = 0.7419354838709677
$ws.Range("D9").Value = 0.7419354838709677
$ws.Range("E9").Value = 0.7419354838709677

$ws.Range("B10").Value = 0.7430555555555556
$ws.Range("C10").Value = 0.7435064935064934
$ws.Range("D10").Value = 0.7419056429232193
$ws.Range("E10").Value = 93

$ws.Range("B11").Value = 0.7449223416965354
$ws.Range("C11").Value = 0.7419354838709677
$ws.Range("D11").Value = 0.7420548476619618
$ws.Range("E11").Value = 93

$ws.Range("B12").Value = 0.9148936170212766
$ws.Range("C12").Value = 0.8775510204081632
$ws.Range("D12").Value = 0.8958333333333333
$ws.Range("E12").Value = 49

$ws.Range("B13").Value = 0.8695652173913043
$ws.Range("C13").Value = 0.9090909090909091
$ws.Range("D13").Value = 0.888888888888889
$ws.Range("E13").Value = 44

$ws.Range("B14").Value = 0.8924731182795699
$ws.Range("C14").Value = 0.8924731182795699
$ws.Range("D14").Value = 0.8924731182795699
$ws.Range("E14").Value = 0.8924731182795699

$ws.Range("B15").Value = 0.8922294172062905
$ws.Range("C15").Value = 0.8933209647495362
$ws.Range("D15").Value = 0.8923611111111112
$ws.Range("E15").Value = 93

$ws.Range("B16").Value = 0.8934479225726876
$ws.Range("C16").Value = 0.8924731182795699
$ws.Range("D16").Value = 0.8925477897252091
$ws.Range("E16").Value = 93

$ws.Range("B17").Value = 0.8409090909090909
$ws.Range("C17").Value = 0.7551020408163265
$ws.Range("D17").Value = 0.7956989247311828
$ws.Range("E17").Value = 49

$ws.Range("B18").Value = 0.7551020408163265
$ws.Range("C18").Value = 0.8409090909090909
$ws.Range("D18").Value = 0.7956989247311828
$ws.Range("E18").Value = 44

$ws.Range("B19").Value = 0.7956989247311828
$ws.Range("C19").Value = 0.7956989247311828
$ws.Range("D19").Value = 0.7956989247311828
$ws.Range("E19").Value = 0.7956989247311828

$ws.Range("B20").Value = 0.7980055658627088
$ws.Range("C20").Value = 0.7980055658627088
$ws.Range("D20").Value = 0.7956989247311828
$ws.Range("E20").Value = 93

$ws.Range("B21").Value = 0.8003122069942346
$ws.Range("C21").Value = 0.7956989247311828
$ws.Range("D21").Value = 0.7956989247311828
$ws.Range("E21").Value = 93

$ws.Range("B22").Value = 0.9411764705882353
$ws.Range("C22").Value = 0.9795918367346939
$ws.Range("D22").Value = 0.96
$ws.Range("E22").Value = 49

$ws.Range("B23").Value = 0.9761904761904762
$ws.Range("C23").Value = 0.9318181818181818
$ws.Range("D23").Value = 0.9534883720930233
$ws.Range("E23").Value = 44

$ws.Range("B24").Value = 0.956989247311828
$ws.Range("C24").Value = 0.956989247311828
$ws.Range("D24").Value = 0.956989247311828
$ws.Range("E24").Value = 0.956989247311828

$ws.Range("B25").Value = 0.9586834733893557
$ws.Range("C25").Value = 0.9557050092764379
$ws.Range("D25").Value = 0.9567441860465116
$ws.Range("E25").Value = 93

$ws.Range("B26").Value = 0.957742236679618
$ws.Range("C26").Value = 0.956989247311828
$ws.Range("D26").Value = 0.9569192298074519
$ws.Range("E26").Value = 93

